$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.389.37"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "1.840.36"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  +1.34%  "
$ws.Range("D5").Value = "'314.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("D7").Value = "'0.4773"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").Value = "'0.3703"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").Value = "'0.07473"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "'0.8862"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "1.852.14"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "'0.07377"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.31%  "
$ws.Range("D14").Value = "'5.488"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "'93.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "'6.591"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "'1.016"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'1.014"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").Value = "'14.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").Value = "27.407.13"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").Value = "'5.353"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").Value = "'10.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").Value = "2.071.13"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").Value = "'1.908"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'152.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("D27").Value = "'18.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "'5.273"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("D30").Value = "'118.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").Value = "'0.08995"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "'0.7607"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("D33").Value = "'1.178"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").Value = "'2.951"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").Value = "'0.05389"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("D41").Value = "'7.316"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").Value = "'0.5361"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").Value = "'2.380"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").Value = "'0.1669"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").Value = "'0.4992"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").Value = "'10.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").Value = "'105.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").Value = "'1.684"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "'0.06323"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
